# Weekly refresh of the "Hortaliza, Feria Lagunitas de Puerto Montt - Albahaca"
# data sheet: a new weekly price observation is inserted as row 136
# (pushing the existing rows 136-174 down to 137-175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136; this shifts rows 136..174 down to 137..175
# and extends the used range from A1:R174 to A1:R175 automatically.
$ws.Rows(136).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(136, 1).Value  = 4
$ws.Cells.Item(136, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(136, 3).Value  = "Los Lagos"
$ws.Cells.Item(136, 4).Value  = 44988
$ws.Cells.Item(136, 5).Value  = 10
$ws.Cells.Item(136, 6).Value  = 100112052
$ws.Cells.Item(136, 7).Value  = "Albahaca"
$ws.Cells.Item(136, 8).Value  = "Sin especificar"
$ws.Cells.Item(136, 9).Value  = "Primera"
$ws.Cells.Item(136, 10).Value = 90
$ws.Cells.Item(136, 11).Value = 6000
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = 6000
$ws.Cells.Item(136, 14).Value = "$/docena de matas"
$ws.Cells.Item(136, 15).Value = "Región Metropolitana"
$ws.Cells.Item(136, 16).Value = 1000
$ws.Cells.Item(136, 17).Value = 6
$ws.Cells.Item(136, 18).Value = "Hortaliza"
